# Apply the "07.Feb.2021" daily-stats update to the COVID Slovakia sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing rows (columns H = AgTests, I = AgPosit) ---
$updates = @{
    310 = @{ H = 74794; I = 3916 }
    311 = @{ H = 62461; I = 1992 }
    320 = @{ I = 3696 }
    323 = @{ H = 149164 }
    324 = @{ H = 231289; I = 2648 }
    325 = @{ H = 705128; I = 5845 }
    326 = @{ H = 417163; I = 3685 }
    327 = @{ H = 235537; I = 2868 }
    328 = @{ H = 178276; I = 2610 }
    329 = @{ H = 82192 }
    331 = @{ H = 147538; I = 2558 }
    332 = @{ H = 411930; I = 4026 }
    333 = @{ H = 252629; I = 2705 }
    334 = @{ H = 201634; I = 3349 }
    335 = @{ H = 121935; I = 2798 }
    336 = @{ H = 96286; I = 3105 }
}

foreach ($rowNum in $updates.Keys) {
    $cols = $updates[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $ws.Range("$colLetter$rowNum").Value = $cols[$colLetter]
    }
}

# --- Append the new row 337 (date 2021-02-04 entry) ---
$ws.Range("A337").Value = 44231
$ws.Range("A337").NumberFormat = "yyyy-mm-dd"
$ws.Range("B337").Value = 259533
$ws.Range("C337").Value = 234371
$ws.Range("D337").Value = 20112
$ws.Range("E337").Value = 11282
$ws.Range("F337").Value = 2630
$ws.Range("G337").Value = 5050
$ws.Range("H337").Value = 93385
$ws.Range("I337").Value = 2979
